$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 408 all contained the same
# date serial 45172 (2023-09-03) and were updated to 45175 (2023-09-06).
$ws.Range("C2:C408").Value = 45175
